$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "SIF 3.0 Framework Installation and Verification" -> "Installation
#    and Verification", with a _GoBack bookmark inserted right before
#    the remaining run. Word enforces unique bookmark names, so adding
#    a new "_GoBack" bookmark here also relocates/removes the one that
#    used to sit right after "... DIRECT environment".
# ------------------------------------------------------------------
$oldHeading = "SIF 3.0 Framework Installation and Verification"

$findRng = $d.Content
$found = $findRng.Find.Execute($oldHeading, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($findRng.Start, $findRng.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
    Write-Output "Inserted _GoBack bookmark before the heading run."
} else {
    Write-Output "WARNING: heading text for _GoBack placement not found."
}

$headingRng = $d.Content
$replaced = $headingRng.Find.Execute($oldHeading, $true, $false, $false, $false, $false, $true, 1, $false, "Installation and Verification", 2)
Write-Output "Heading text replaced: $replaced"

# ------------------------------------------------------------------
# 2) Merge the three runs of the "Scheduling SIF Events" paragraph into
#    a single run (identical visible text, no mid-sentence run splits).
# ------------------------------------------------------------------
$schedText = "As the implementation of SIF Events in the SIF Framework is done by using a REST endpoint, it is the responsibility of the developer to schedule calls to execute the broadcast of change records. The mechanism used for scheduling SIF Events is therefore outside the scope of this document."
$schedRng = $d.Content
$schedReplaced = $schedRng.Find.Execute($schedText, $true, $false, $false, $false, $false, $true, 1, $false, $schedText, 2)
Write-Output "Scheduling SIF Events paragraph runs merged: $schedReplaced"

# ------------------------------------------------------------------
# 3) Bump the cached "SystemVersion" field results in the odd/even page
#    footers from 3.2.1.3 to 3.2.1.4 (the title-page DOCPROPERTY field
#    on the cover page is left untouched). Editing the field's Result
#    range directly (instead of Find/Replace across the whole footer)
#    avoids Word re-normalising unrelated runs/fields in that footer.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        foreach ($fld in $ftr.Range.Fields) {
            if ($fld.Code.Text -match "SystemVersion") {
                $fld.Result.Text = "3.2.1.4"
                Write-Output "Updated SystemVersion field result in footer $i to 3.2.1.4"
            }
        }
    }
}
